$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Enhancements", 1),
    @("Hud", 1),
    @("Follow/Remain Button", 0.5),
    @("Persuading an enemy to follow you", 1),
    @("Enemy shoot back when syndicate's gun out and are in tange", 1),
    @("Enemy drop random items you pick up by walking over them", 0.5),
    @("Spawn enemy's at random position", "15min")
)

$row = 11
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $row++
}

$ws.Range("A18").Select()
